$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detektory")

# Column G (index 7) becomes 20 for every data row (2-76)
$ws.Range("G2:G76").Value = 20

# Column H (index 8) values per row, taken from the target diff
$hValues = @{
    2 = 1
    3 = 1
    4 = 0
    5 = 0
    6 = 1
    7 = 1
    8 = 0
    9 = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 1
    59 = 1
    60 = 0
    61 = 0
    62 = 1
    63 = 1
    64 = 0
    65 = 0
    66 = 0
    67 = 1
    68 = 1
    69 = 0
    70 = 0
    71 = 1
    72 = 1
    73 = 0
    74 = 0
    75 = 1
    76 = 1
}

foreach ($row in $hValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $hValues[$row]
}

# Update the selected cell in the sheet view to H16
$ws.Range("H16").Select()
